# Updates odds/snapshot data in the active worksheet to reflect the latest
# Betfair Back/Lay snapshot refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 3.35
$ws.Range("I2").Value = 2.22
$ws.Range("J2").Value = 2.76
$ws.Range("Q2").Value = 2.12
$ws.Range("BH2").Value = "2026-02-22 08:31:36"

# Row 3
$ws.Range("H3").Value = 4.4
$ws.Range("P3").Value = 1.93
$ws.Range("Q3").Value = 1.86
$ws.Range("BH3").Value = "2026-02-22 08:31:36"

# Row 4
$ws.Range("G4").Value = 1.51
$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 5.1
$ws.Range("P4").Value = 2.82
$ws.Range("Q4").Value = 1.4
$ws.Range("BH4").Value = "2026-02-22 08:31:36"

# Row 5
$ws.Range("F5").Value = 3.4
$ws.Range("G5").Value = 4.3
$ws.Range("H5").Value = 1.97
$ws.Range("BH5").Value = "2026-02-22 08:31:36"

# Row 6
$ws.Range("G6").Value = 11
$ws.Range("I6").Value = 1.45
$ws.Range("K6").Value = 6.2
$ws.Range("BH6").Value = "2026-02-22 08:31:36"

# Row 7
$ws.Range("F7").Value = 1.73
$ws.Range("G7").Value = 2.18
$ws.Range("H7").Value = 3.5
$ws.Range("J7").Value = 3.65
$ws.Range("K7").Value = 5.5
$ws.Range("P7").Value = 2.04
$ws.Range("Q7").Value = 1.57
$ws.Range("BH7").Value = "2026-02-22 08:31:36"

# Row 8
$ws.Range("G8").Value = 9.199999999999999
$ws.Range("H8").Value = 1.46
$ws.Range("I8").Value = 1.58
$ws.Range("P8").Value = 2.28
$ws.Range("Q8").Value = 1.61
$ws.Range("BH8").Value = "2026-02-22 08:31:36"

# Row 9
$ws.Range("I9").Value = 1.6
$ws.Range("BH9").Value = "2026-02-22 08:31:36"

# Row 10
$ws.Range("F10").Value = 1.53
$ws.Range("G10").Value = 1.63
$ws.Range("J10").Value = 3.75
$ws.Range("K10").Value = 4.4
$ws.Range("Q10").Value = 2.46
$ws.Range("BH10").Value = "2026-02-22 08:31:36"

# Row 11
$ws.Range("F11").Value = 3.65
$ws.Range("H11").Value = 2.36
$ws.Range("I11").Value = 2.76
$ws.Range("J11").Value = 2.76
$ws.Range("K11").Value = 3.35
$ws.Range("P11").Value = 1.42
$ws.Range("Q11").Value = 2.88
$ws.Range("BH11").Value = "2026-02-22 08:31:36"

# Row 12
$ws.Range("I12").Value = 4.5
$ws.Range("K12").Value = 3.15
$ws.Range("BH12").Value = "2026-02-22 08:31:36"

# Row 13
$ws.Range("F13").Value = 1.44
$ws.Range("H13").Value = 1.44
$ws.Range("K13").Value = 3.3
$ws.Range("BH13").Value = "2026-02-22 08:31:36"

# Row 14 (timestamp only)
$ws.Range("BH14").Value = "2026-02-22 08:31:36"

# Row 15
$ws.Range("I15").Value = 2.68
$ws.Range("K15").Value = 3.2
$ws.Range("BH15").Value = "2026-02-22 08:31:36"
